{"js": "// ONC-1272: Restore ability to control authoring and group access permissions\n// from Assignments Beta.\n//\n// Insert a new \"GRADER PERMISSION SETTINGS (GPS-1)\" row into the TOC table,\n// immediately before the existing \"GRADER PERMISSIONS HELPER (GPH-1)\" row,\n// with page number 23. Because the new entry pushes every following entry\n// forward by two pages, every page number from that point on (inclusive of\n// the former GPH-1 row, now on page 25) is incremented by 2.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Locate the \"GRADER PERMISSIONS HELPER (GPH-1)\" row.\nlet gphIndex = -1;\nfor (let i = 0; i < table.values.length; i++) {\n  if (table.values[i][0].indexOf(\"GRADER PERMISSIONS HELPER\") !== -1) {\n    gphIndex = i;\n    break;\n  }\n}\nif (gphIndex === -1) {\n  throw new Error(\"Could not find GRADER PERMISSIONS HELPER row\");\n}\n\n// Every row at/after the GPH-1 row gets its page number bumped by 2.\nfor (let i = gphIndex; i < table.values.length; i++) {\n  const oldPage = parseInt(table.values[i][1], 10);\n  const newPage = (oldPage + 2).toString();\n  table.getCell(i, 1).value = newPage;\n}\nawait context.sync();\n\n// Insert the new row immediately above the GPH-1 row.\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst gphRow = table.rows.items[gphIndex];\ngphRow.insertRows(\"Before\", 1, [[\"GRADER PERMISSION SETTINGS (GPS-1)\", \"23\"]]);\nawait context.sync();\n", "ps1": "# ONC-1272: Restore ability to control authoring and group access permissions\n# from Assignments Beta.\n#\n# Insert a new \"GRADER PERMISSION SETTINGS (GPS-1)\" row into the TOC table,\n# immediately before the existing \"GRADER PERMISSIONS HELPER (GPH-1)\" row,\n# with page number 23. Because the new entry pushes every following entry\n# forward by two pages, every page number from that point on (inclusive of\n# the former GPH-1 row, now pushed to page 25) is incremented by 2.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the \"GRADER PERMISSIONS HELPER (GPH-1)\" row by scanning cell text.\n$gphRowIndex = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $rowText = $t.Cell($i, 1).Range.Text\n    if ($rowText -like \"*GRADER PERMISSIONS HELPER*\") {\n        $gphRowIndex = $i\n        break\n    }\n}\n\nif ($gphRowIndex -eq -1) {\n    throw \"Could not find GRADER PERMISSIONS HELPER row\"\n}\n\n# Every row at/after the GPH-1 row gets its page number bumped by 2 -- do\n# this first (from the bottom up) before inserting, so indices stay valid.\nfor ($i = $t.Rows.Count; $i -ge $gphRowIndex; $i--) {\n    $pageCell = $t.Cell($i, 2)\n    $oldPage = [int]($pageCell.Range.Text -replace \"[^0-9]\", \"\")\n    $newPage = $oldPage + 2\n    $pageCell.Range.Text = [string]$newPage\n}\n\n# Insert the new row immediately above the GPH-1 row (which has now shifted\n# down by one row due to the loop above not touching row count).\n$newRow = $t.Rows.Add($t.Rows.Item($gphRowIndex))\n$newRow.Cells.Item(1).Range.Text = \"GRADER PERMISSION SETTINGS (GPS-1)\"\n$newRow.Cells.Item(2).Range.Text = \"23\"\n"}
